# Apply updated Fruta/Hortaliza weekly price data (Vega Modelo de Temuco - Caigua)
# Updates columns D (Fecha), J (Volumen), K/L/M (Precio min/max/promedio), P (Precio $/Kg)
# for rows 2-28 of the active worksheet, per the new source extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, J2, K2, L2, M2, P2
$ws.Cells.Item(2, 4).Value = 44769
$ws.Cells.Item(2, 10).Value = 50
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 16).Value = 1333

# Row 3: D3, J3, K3, L3, M3, P3
$ws.Cells.Item(3, 4).Value = 44845
$ws.Cells.Item(3, 10).Value = 20
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 16000
$ws.Cells.Item(3, 16).Value = 1067

# Row 4: D4
$ws.Cells.Item(4, 4).Value = 44826

# Row 5: D5, J5
$ws.Cells.Item(5, 4).Value = 44819
$ws.Cells.Item(5, 10).Value = 100

# Row 6: D6
$ws.Cells.Item(6, 4).Value = 44813

# Row 7: D7, J7, K7, L7, M7, P7
$ws.Cells.Item(7, 4).Value = 44841
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 16000
$ws.Cells.Item(7, 16).Value = 1067

# Row 8: D8, J8, K8, L8, M8, P8
$ws.Cells.Item(8, 4).Value = 44830
$ws.Cells.Item(8, 10).Value = 25
$ws.Cells.Item(8, 11).Value = 12000
$ws.Cells.Item(8, 12).Value = 12000
$ws.Cells.Item(8, 13).Value = 12000
$ws.Cells.Item(8, 16).Value = 800

# Row 9: D9, J9
$ws.Cells.Item(9, 4).Value = 44838
$ws.Cells.Item(9, 10).Value = 10

# Row 10: D10, K10, L10, M10, P10
$ws.Cells.Item(10, 4).Value = 44755
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 20000
$ws.Cells.Item(10, 16).Value = 1333

# Row 11: D11, J11, K11, L11, M11, P11
$ws.Cells.Item(11, 4).Value = 44508
$ws.Cells.Item(11, 10).Value = 40
$ws.Cells.Item(11, 11).Value = 10000
$ws.Cells.Item(11, 12).Value = 10000
$ws.Cells.Item(11, 13).Value = 10000
$ws.Cells.Item(11, 16).Value = 667

# Row 12: D12, J12, K12, L12, M12, P12
$ws.Cells.Item(12, 4).Value = 44525
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 8000
$ws.Cells.Item(12, 12).Value = 8000
$ws.Cells.Item(12, 13).Value = 8000
$ws.Cells.Item(12, 16).Value = 533

# Row 13: D13, J13
$ws.Cells.Item(13, 4).Value = 44827
$ws.Cells.Item(13, 10).Value = 20

# Row 14: D14, J14
$ws.Cells.Item(14, 4).Value = 44771
$ws.Cells.Item(14, 10).Value = 40

# Row 15: D15, J15, K15, L15, M15, P15
$ws.Cells.Item(15, 4).Value = 44518
$ws.Cells.Item(15, 10).Value = 50
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 16).Value = 667

# Row 16: D16, J16, K16, L16, M16, P16
$ws.Cells.Item(16, 4).Value = 44839
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 16000
$ws.Cells.Item(16, 16).Value = 1067

# Row 17: D17, J17
$ws.Cells.Item(17, 4).Value = 44749
$ws.Cells.Item(17, 10).Value = 50

# Row 18: D18, J18, K18, L18, M18, P18
$ws.Cells.Item(18, 4).Value = 45134
$ws.Cells.Item(18, 10).Value = 5
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 20000
$ws.Cells.Item(18, 13).Value = 20000
$ws.Cells.Item(18, 16).Value = 1333

# Row 19: D19, J19
$ws.Cells.Item(19, 4).Value = 44812
$ws.Cells.Item(19, 10).Value = 80

# Row 20: D20, J20, K20, L20, M20, P20
$ws.Cells.Item(20, 4).Value = 44756
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 16).Value = 1333

# Row 21: D21, J21
$ws.Cells.Item(21, 4).Value = 45225
$ws.Cells.Item(21, 10).Value = 80

# Row 22: D22, J22, K22, L22, M22, P22
$ws.Cells.Item(22, 4).Value = 44837
$ws.Cells.Item(22, 10).Value = 80
$ws.Cells.Item(22, 11).Value = 16000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 16000
$ws.Cells.Item(22, 16).Value = 1067

# Row 23: D23, J23
$ws.Cells.Item(23, 4).Value = 44811
$ws.Cells.Item(23, 10).Value = 30

# Row 24: D24, J24, K24, L24, M24, P24
$ws.Cells.Item(24, 4).Value = 44757
$ws.Cells.Item(24, 10).Value = 30
$ws.Cells.Item(24, 11).Value = 20000
$ws.Cells.Item(24, 12).Value = 20000
$ws.Cells.Item(24, 13).Value = 20000
$ws.Cells.Item(24, 16).Value = 1333

# Row 25: D25, J25, K25, L25, M25, P25
$ws.Cells.Item(25, 4).Value = 44767
$ws.Cells.Item(25, 10).Value = 50
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 13).Value = 20000
$ws.Cells.Item(25, 16).Value = 1333

# Row 26: D26, J26
$ws.Cells.Item(26, 4).Value = 44825
$ws.Cells.Item(26, 10).Value = 30

# Row 27: D27, J27
$ws.Cells.Item(27, 4).Value = 44776
$ws.Cells.Item(27, 10).Value = 80

# Row 28: D28, J28, K28, L28, M28, P28
$ws.Cells.Item(28, 4).Value = 44824
$ws.Cells.Item(28, 10).Value = 20
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 20000
$ws.Cells.Item(28, 13).Value = 20000
$ws.Cells.Item(28, 16).Value = 1333
